# Auto-generated Excel COM-interop script
# Applies the "Updated symbol list" diff: refreshed Price/Volume(1h)/Hora
# columns for the cryptos.xlsx coin table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '272.38'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.71%'
$ws.Range("E2").Style = "Normal"

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '22'
$ws.Range("G2").Style = "Normal"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '0.32%'
$ws.Range("E3").Style = "Normal"

$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '22'
$ws.Range("G3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.910'
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '4.13%'
$ws.Range("E4").Style = "Normal"

$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '22'
$ws.Range("G4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06310'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '3.10%'
$ws.Range("E5").Style = "Normal"

$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '22'
$ws.Range("G5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.892'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '2.22%'
$ws.Range("E6").Style = "Normal"

$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '22'
$ws.Range("G6").Style = "Normal"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '5.40%'
$ws.Range("E7").Style = "Normal"

$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '22'
$ws.Range("G7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.401'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '56.76%'
$ws.Range("E8").Style = "Normal"

$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '22'
$ws.Range("G8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8840'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '3.31%'
$ws.Range("E9").Style = "Normal"

$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '22'
$ws.Range("G9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1469'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2.50%'
$ws.Range("E10").Style = "Normal"

$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '22'
$ws.Range("G10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.05076'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '0.27%'
$ws.Range("E11").Style = "Normal"

$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '22'
$ws.Range("G11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07376'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '3.62%'
$ws.Range("E12").Style = "Normal"

$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '22'
$ws.Range("G12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03169'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.51%'
$ws.Range("E13").Style = "Normal"

$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '22'
$ws.Range("G13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09037'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.00%'
$ws.Range("E14").Style = "Normal"

$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '22'
$ws.Range("G14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001583'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '3.54%'
$ws.Range("E15").Style = "Normal"

$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '22'
$ws.Range("G15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0006332'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '4.58%'
$ws.Range("E16").Style = "Normal"

$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '22'
$ws.Range("G16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.006064'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.69%'
$ws.Range("E17").Style = "Normal"

$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '22'
$ws.Range("G17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.471'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.22%'
$ws.Range("E18").Style = "Normal"

$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '22'
$ws.Range("G18").Style = "Normal"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.89%'
$ws.Range("E19").Style = "Normal"

$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '22'
$ws.Range("G19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3143'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '1.71%'
$ws.Range("E20").Style = "Normal"

$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '22'
$ws.Range("G20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1335'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '4.19%'
$ws.Range("E21").Style = "Normal"

$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '22'
$ws.Range("G21").Style = "Normal"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.90%'
$ws.Range("E22").Style = "Normal"

$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '22'
$ws.Range("G22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04340'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '2.16%'
$ws.Range("E23").Style = "Normal"

$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '22'
$ws.Range("G23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001178'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.03%'
$ws.Range("E24").Style = "Normal"

$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '22'
$ws.Range("G24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.003637'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-12.28%'
$ws.Range("E25").Style = "Normal"

$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '22'
$ws.Range("G25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001202'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.05%'
$ws.Range("E26").Style = "Normal"

$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '22'
$ws.Range("G26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001700'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '1.08%'
$ws.Range("E27").Style = "Normal"

$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '22'
$ws.Range("G27").Style = "Normal"

$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '22'
$ws.Range("G28").Style = "Normal"

$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '22'
$ws.Range("G29").Style = "Normal"

$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '22'
$ws.Range("G30").Style = "Normal"

$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '22'
$ws.Range("G31").Style = "Normal"

$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '22'
$ws.Range("G32").Style = "Normal"

$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '22'
$ws.Range("G33").Style = "Normal"

$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '22'
$ws.Range("G34").Style = "Normal"

$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '22'
$ws.Range("G35").Style = "Normal"

$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '22'
$ws.Range("G36").Style = "Normal"

$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '22'
$ws.Range("G37").Style = "Normal"

$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '22'
$ws.Range("G38").Style = "Normal"

$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '22'
$ws.Range("G39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04049'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2.53%'
$ws.Range("E40").Style = "Normal"

$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '22'
$ws.Range("G40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006617'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '57.69%'
$ws.Range("E41").Style = "Normal"

$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '22'
$ws.Range("G41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1161'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '3.77%'
$ws.Range("E42").Style = "Normal"

$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '22'
$ws.Range("G42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002143'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '5.13%'
$ws.Range("E43").Style = "Normal"

$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '22'
$ws.Range("G43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01259'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '7.94%'
$ws.Range("E44").Style = "Normal"

$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '22'
$ws.Range("G44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005337'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '4.08%'
$ws.Range("E45").Style = "Normal"

$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '22'
$ws.Range("G45").Style = "Normal"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '152.74%'
$ws.Range("E46").Style = "Normal"

$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '22'
$ws.Range("G46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.02124'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-29.06%'
$ws.Range("E47").Style = "Normal"

$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '22'
$ws.Range("G47").Style = "Normal"

$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '22'
$ws.Range("G48").Style = "Normal"

$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '22'
$ws.Range("G49").Style = "Normal"

$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '22'
$ws.Range("G50").Style = "Normal"

$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '22'
$ws.Range("G51").Style = "Normal"
